$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.142.27'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '1.668.25'
$ws.Range("E3").Value = '  -1.61%  '
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.86'
$ws.Range("E5").Value = '  -3.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5244'
$ws.Range("E6").Value = '  -3.63%  '
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2631'
$ws.Range("E8").Value = '  -3.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06301'
$ws.Range("E9").Value = '  -2.61%  '
$ws.Range("E10").Value = '  -2.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07561'
$ws.Range("E11").Value = '  -1.34%  '
$ws.Range("D12").Value = '1.671.13'
$ws.Range("E12").Value = '  -1.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.442'
$ws.Range("E13").Value = '  -2.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5577'
$ws.Range("E14").Value = '  -4.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.99'
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000007934'
$ws.Range("E16").Value = '  -5.94%  '
$ws.Range("D17").Value = '26.168.04'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.744'
$ws.Range("E19").Value = '  -3.69%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '186.83'
$ws.Range("E20").Value = '  -2.51%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.36'
$ws.Range("E21").Value = '  -5.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.179'
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.005'
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.80'
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1248'
$ws.Range("E25").Value = '  -3.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.531'
$ws.Range("E26").Value = '  -4.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.97'
$ws.Range("E27").Value = '  +0.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06245'
$ws.Range("E28").Value = '  -2.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.342'
$ws.Range("E29").Value = '  -3.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.281'
$ws.Range("E30").Value = '  -3.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.507'
$ws.Range("E31").Value = '  -2.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.414'
$ws.Range("E32").Value = '  -4.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.628'
$ws.Range("E33").Value = '  -3.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9956'
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6040'
$ws.Range("E35").Value = '  -2.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.409'
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.737'
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.109.51'
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.117'
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01609'
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8765'
$ws.Range("E41").Value = '  -1.00%  '
$ws.Range("E42").Value = '  -0.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.90'
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("D44").Value = '1.823.50'
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000111'
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.41'
$ws.Range("E46").Value = '  -4.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9995'
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.012'
$ws.Range("E48").Value = '  -2.47%  '
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4251'
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.977'
$ws.Range("E51").Value = '  -2.09%  '
